$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.496418714523315
$ws.Range("B1").Value = 6.587157249450684
$ws.Range("C1").Value = 2.591421842575073
$ws.Range("D1").Value = 1.696979284286499
$ws.Range("E1").Value = 1.357116103172302
